$wb = $excel.ActiveWorkbook

# --- Sheet 1: ID_7e0f17b -------------------------------------------------
# Fix the date on the existing row, then duplicate it into row 3.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ID_7e0f17b"
$ws1.Range("B2").Value = "15/02/2025"
$ws1.Range("A2:E2").Copy() | Out-Null
$ws1.Range("A3").PasteSpecial() | Out-Null

# --- Sheet 2: ID_dda333c -------------------------------------------------
# Duplicate sheet 1 (keeps header formatting), rename, and replace data rows.
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "ID_dda333c"

$ws2.Range("A2").Value = 22.3
$ws2.Range("B2").Value = "15/02/2025"
$ws2.Range("C2").Value = "Mens Casual Premium Slim Fit T-Shirts "
$ws2.Range("D2").Value = "Slim-fitting style, contrast raglan long sleeve, three-button henley placket, light weight & soft fabric for breathable and comfortable wearing. And Solid stitched shirts with round neck made for durability and a great fit for casual fashion wear and diehard baseball fans. The Henley style round neckline includes a three-button placket."
$ws2.Range("E2").Value = "ID_dda333c"

$ws2.Range("A3").Value = 22.3
$ws2.Range("B3").Value = "15/02/2025"
$ws2.Range("C3").Value = "Mens Casual Premium Slim Fit T-Shirts "
$ws2.Range("D3").Value = "Slim-fitting style, contrast raglan long sleeve, three-button henley placket, light weight & soft fabric for breathable and comfortable wearing. And Solid stitched shirts with round neck made for durability and a great fit for casual fashion wear and diehard baseball fans. The Henley style round neckline includes a three-button placket."
$ws2.Range("E3").Value = "ID_dda333c"

# --- Sheet 3: ID_ca64249 --------------------------------------------------
# A single-row sheet (no row 3), replace the product info.
$ws2.Copy($null, $ws2) | Out-Null
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "ID_ca64249"

$ws3.Range("A3:E3").ClearContents() | Out-Null

$ws3.Range("A2").Value = 109
$ws3.Range("B2").Value = "15/02/2025"
$ws3.Range("C2").Value = "SanDisk SSD PLUS 1TB Internal SSD - SATA III 6 Gb/s"
$ws3.Range("D2").Value = "Easy upgrade for faster boot up, shutdown, application load and response (As compared to 5400 RPM SATA 2.5” hard drive; Based on published specifications and internal benchmarking tests using PCMark vantage scores) Boosts burst write performance, making it ideal for typical PC workloads The perfect balance of performance and reliability Read/write speeds of up to 535MB/s/450MB/s (Based on internal testing; Performance may vary depending upon drive capacity, host device, OS and application.)"
$ws3.Range("E2").Value = "ID_ca64249"
